{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the five wording/technology updates described by the diff:\n//   1. \"Proficient in PHP, SQL, MVC, Python\"\n//        -> \"Proficient in React, PHP, SQL, MVC, Python\"\n//   2. \"Modernizing UI platform for OneID system using React, Redux, Node.js\"\n//        -> \"...React, Redux, Node.js, Jest\"\n//   3. \"Developing and automating CI CD for various projects using Gitlab, Docker\"\n//        -> \"...Gitlab, Docker, Bash Scripts\"\n//   4. \"Used Fireworks, Photoshop, and Gimp to create UI elements from mockups\"\n//        -> \"Created custom UI elements from mockups in Fireworks, Photoshop, and Gimp\"\n//   5. \"Oracle PL/SQL client/server final project and an Amazon.com like clone\n//       for a final software engineering project using PHP\"\n//        -> \"Oracle PL/SQL client-server final project and PHP backed Amazon.com\n//            like e-commerce final software engineering project\"\n//\n// (The diff also shows the two header icon images (envelope.png / github.png)\n// swapping their internal `name` attribute (\"image2.png\" <-> \"image3.png\") in\n// <wp:docPr>/<pic:cNvPr>. That attribute is purely an internal display name\n// for the drawing and is not exposed by the Word JavaScript API --\n// Word.InlinePicture has no settable/gettable \"name\" property in Office.js\n// (only altTextDescription / altTextTitle are exposed), so that cosmetic\n// rename cannot be performed through this API and is intentionally left\n// untouched here.)\n\nconst replacements = [\n  {\n    find: \"Proficient in PHP, SQL, MVC, Python\",\n    replace: \"Proficient in React, PHP, SQL, MVC, Python\",\n  },\n  {\n    find: \"Modernizing UI platform for OneID system using React, Redux, Node.js\",\n    replace: \"Modernizing UI platform for OneID system using React, Redux, Node.js, Jest\",\n  },\n  {\n    find: \"Developing and automating CI CD for various projects using Gitlab, Docker\",\n    replace: \"Developing and automating CI CD for various projects using Gitlab, Docker, Bash Scripts\",\n  },\n  {\n    find: \"Used Fireworks, Photoshop, and Gimp to create UI elements from mockups\",\n    replace: \"Created custom UI elements from mockups in Fireworks, Photoshop, and Gimp\",\n  },\n  {\n    find: \"Oracle PL/SQL client/server final project and an Amazon.com like clone for a final software engineering project using PHP\",\n    replace: \"Oracle PL/SQL client-server final project and PHP backed Amazon.com like e-commerce final software engineering project\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found, cannot apply replacement: \" + find);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the five wording/technology updates described by the diff:\n#   1. \"Proficient in PHP, SQL, MVC, Python\"\n#        -> \"Proficient in React, PHP, SQL, MVC, Python\"\n#   2. \"Modernizing UI platform for OneID system using React, Redux, Node.js\"\n#        -> \"...React, Redux, Node.js, Jest\"\n#   3. \"Developing and automating CI CD for various projects using Gitlab, Docker\"\n#        -> \"...Gitlab, Docker, Bash Scripts\"\n#   4. \"Used Fireworks, Photoshop, and Gimp to create UI elements from mockups\"\n#        -> \"Created custom UI elements from mockups in Fireworks, Photoshop, and Gimp\"\n#   5. \"Oracle PL/SQL client/server final project and an Amazon.com like clone\n#       for a final software engineering project using PHP\"\n#        -> \"Oracle PL/SQL client-server final project and PHP backed Amazon.com\n#            like e-commerce final software engineering project\"\n#\n# (The diff also shows the two header icon images (envelope.png / github.png)\n# swapping their internal `name` attribute (\"image2.png\" <-> \"image3.png\") in\n# <wp:docPr>/<pic:cNvPr>. That attribute is purely an internal display name\n# for the drawing and is not exposed by the Word object model -- InlineShape\n# (unlike the floating Shape object) has no Name property in Word's COM/VBA\n# API, only .Title/.AlternativeText, so that cosmetic rename cannot be\n# performed through this API and is intentionally left untouched here.)\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindWrapNone = 0\n\n$replacements = @(\n    @{\n        Find    = \"Proficient in PHP, SQL, MVC, Python\"\n        Replace = \"Proficient in React, PHP, SQL, MVC, Python\"\n    },\n    @{\n        Find    = \"Modernizing UI platform for OneID system using React, Redux, Node.js\"\n        Replace = \"Modernizing UI platform for OneID system using React, Redux, Node.js, Jest\"\n    },\n    @{\n        Find    = \"Developing and automating CI CD for various projects using Gitlab, Docker\"\n        Replace = \"Developing and automating CI CD for various projects using Gitlab, Docker, Bash Scripts\"\n    },\n    @{\n        Find    = \"Used Fireworks, Photoshop, and Gimp to create UI elements from mockups\"\n        Replace = \"Created custom UI elements from mockups in Fireworks, Photoshop, and Gimp\"\n    },\n    @{\n        Find    = \"Oracle PL/SQL client/server final project and an Amazon.com like clone for a final software engineering project using PHP\"\n        Replace = \"Oracle PL/SQL client-server final project and PHP backed Amazon.com like e-commerce final software engineering project\"\n    }\n)\n\nforeach ($item in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $item.Find\n    $find.Replacement.Text = $item.Replace\n\n    $ok = $find.Execute($item.Find, $false, $true, $false, $false, $false, $true, $wdFindWrapNone, $false, $item.Replace, $wdReplaceAll)\n    if (-not $ok) {\n        throw \"Text not found, cannot apply replacement: $($item.Find)\"\n    }\n}\n"}
